# Updated cryptos list (Price + Volume(1h) columns) per upstream data refresh.
# Each assignment is prefixed with a literal apostrophe so Excel stores the
# new reading as text (matching the sheet's existing inlineStr/string cells)
# instead of silently re-typing it as a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'34.709.32"
$ws.Range('D3').Value = "'1.809.64"
$ws.Range('E3').Value = "'  -1.78%  "
$ws.Range('E4').Value = "'  +0.29%  "
$ws.Range('D5').Value = "'232.43"
$ws.Range('E5').Value = "'  +1.24%  "
$ws.Range('D6').Value = "'0.604"
$ws.Range('E6').Value = "'  -0.93%  "
$ws.Range('E7').Value = "'  +0.30%  "
$ws.Range('D8').Value = "'39.28"
$ws.Range('E8').Value = "'  -9.32%  "
$ws.Range('D9').Value = "'0.326"
$ws.Range('E9').Value = "'  +5.82%  "
$ws.Range('D11').Value = "'0.0991"
$ws.Range('E11').Value = "'  -1.96%  "
$ws.Range('D12').Value = "'2.070.24"
$ws.Range('E12').Value = "'  -1.85%  "
$ws.Range('E13').Value = "'  +0.00%  "
$ws.Range('D14').Value = "'1.821.51"
$ws.Range('E14').Value = "'  -1.09%  "
$ws.Range('D15').Value = "'11.08"
$ws.Range('E15').Value = "'  -2.01%  "
$ws.Range('E16').Value = "'  -1.90%  "
$ws.Range('D17').Value = "'34.685.95"
$ws.Range('E17').Value = "'  -2.10%  "
$ws.Range('D18').Value = "'69.45"
$ws.Range('E18').Value = "'  -0.98%  "
$ws.Range('D19').Value = "'0.0₃0786"
$ws.Range('E19').Value = "'  -1.05%  "
$ws.Range('D20').Value = "'240.02"
$ws.Range('E20').Value = "'  -1.84%  "
$ws.Range('D21').Value = "'11.97"
$ws.Range('E21').Value = "'  -1.19%  "
$ws.Range('D22').Value = "'4.68"
$ws.Range('E22').Value = "'  +0.37%  "
$ws.Range('E23').Value = "'  +0.36%  "
$ws.Range('E24').Value = "'  +2.17%  "
$ws.Range('D25').Value = "'172.18"
$ws.Range('E25').Value = "'  +1.99%  "
$ws.Range('E26').Value = "'  -2.47%  "
$ws.Range('D27').Value = "'17.18"
$ws.Range('E27').Value = "'  -3.34%  "
$ws.Range('D28').Value = "'0.120"
$ws.Range('E28').Value = "'  -1.78%  "
$ws.Range('E29').Value = "'  +11.97%  "
$ws.Range('E30').Value = "'  +0.34%  "
$ws.Range('D31').Value = "'4.03"
$ws.Range('E31').Value = "'  +2.77%  "
$ws.Range('E32').Value = "'  +0.03%  "
$ws.Range('D33').Value = "'3.96"
$ws.Range('E33').Value = "'  -2.65%  "
$ws.Range('D34').Value = "'1.27"
$ws.Range('E34').Value = "'  +17.19%  "
$ws.Range('E35').Value = "'  -4.59%  "
$ws.Range('D36').Value = "'0.699"
$ws.Range('E36').Value = "'  +1.26%  "
$ws.Range('D37').Value = "'91.49"
$ws.Range('E37').Value = "'  -4.48%  "
$ws.Range('E38').Value = "'  +5.37%  "
$ws.Range('D39').Value = "'1.323.12"
$ws.Range('E39').Value = "'  -1.66%  "
$ws.Range('E40').Value = "'  -1.15%  "
$ws.Range('D41').Value = "'2.47"
$ws.Range('E41').Value = "'  +0.67%  "
$ws.Range('E42').Value = "'  -4.24%  "
$ws.Range('D43').Value = "'14.16"
$ws.Range('E43').Value = "'  -5.71%  "
$ws.Range('E44').Value = "'  -9.07%  "
$ws.Range('D45').Value = "'2.67"
$ws.Range('E45').Value = "'  -4.88%  "
$ws.Range('D46').Value = "'6.27"
$ws.Range('E46').Value = "'  +1.01%  "
$ws.Range('E47').Value = "'  -1.33%  "
$ws.Range('D48').Value = "'1.996.93"
$ws.Range('E48').Value = "'  -0.65%  "
$ws.Range('E49').Value = "'  +0.33%  "
$ws.Range('D50').Value = "'0.0668"
$ws.Range('E50').Value = "'  +7.32%  "
$ws.Range('D51').Value = "'98.65"
$ws.Range('E51').Value = "'  -4.13%  "

Write-Host "Updated $($wb.ActiveSheet.Name): 80 cells across D (Price) and E (Volume 1h) columns, rows 2-51"
